$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) to English snake_case field names ---
$ws.Cells.Item(1, 1).Value2 = "mx_state"
$ws.Cells.Item(1, 2).Value2 = "mx_municipality"
$ws.Cells.Item(1, 3).Value2 = "n_matriculas"
$ws.Cells.Item(1, 4).Value2 = "pct_matriculas"

# --- 2. Title-case the small Spanish connector words ("de", "del", "la", ---
#        "las", "los", "el", "y") wherever they appear as a standalone    ---
#        word in columns A (state) and B (municipality), rows 2..1345.   ---
$smallWords = @('de', 'del', 'la', 'las', 'los', 'el', 'y')

function Convert-SmallWords($text) {
    $words = $text.Split(' ')
    $newWords = @()
    foreach ($w in $words) {
        if ($smallWords -contains $w) {
            $newWords += $w.Substring(0, 1).ToUpper() + $w.Substring(1)
        } else {
            $newWords += $w
        }
    }
    return ($newWords -join ' ')
}

for ($r = 2; $r -le 1345; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null) {
        $aCell.Value2 = Convert-SmallWords($aVal)
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null) {
        $bCell.Value2 = Convert-SmallWords($bVal)
    }
}

# --- 3. Grand-total label: "TOTAL" -> "Total" ---
$ws.Cells.Item(1345, 1).Value2 = "Total"

# --- 4. Drop the trailing footnote rows (sample size / source / credits / date) ---
$ws.Rows("1347:1351").Delete()
